$d = $word.ActiveDocument

$d.Content.Find.Execute("Информационный технологии", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Информатика", 2)
